$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.545.63"
$ws.Range("E2").Value = "  +3.31%  "
$ws.Range("D3").Value = "2.256.88"
$ws.Range("E3").Value = "  +1.97%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "91.23"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.67%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.528"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.64%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.478"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.01"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.82"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.92%  "
$ws.Range("E12").Value = "  +1.84%  "
$ws.Range("E13").Value = "  +0.80%  "
$ws.Range("E14").Value = "  +2.13%  "
$ws.Range("D15").Value = "2.602.82"
$ws.Range("E15").Value = "  +1.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.17"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.25%  "
$ws.Range("D17").Value = "2.271.23"
$ws.Range("E17").Value = "  +3.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.758"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.79%  "
$ws.Range("D19").Value = "41.484.79"
$ws.Range("E19").Value = "  +3.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.93%  "
$ws.Range("D21").Value = "0.0₃0901"
$ws.Range("E21").Value = "  +1.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.88"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.90%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "239.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.57"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.52%  "
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.89"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.77%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.06"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.80%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.47"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "160.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "34.06"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.92%  "
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.16"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0739"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("E37").Value = "  +2.29%  "
$ws.Range("B38").Value = "Celestia"
$ws.Range("C38").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "16.61"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.16%  "
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.116"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.104"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.99%  "
$ws.Range("E41").Value = "  +2.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.93"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.77%  "
$ws.Range("D43").Value = "2.045.08"
$ws.Range("E43").Value = "  -1.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.33"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.38%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0277"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.29"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.12%  "
$ws.Range("E47").Value = "  +6.99%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.85"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.57%  "
$ws.Range("E49").Value = "  +3.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "72.42"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.59%  "
